$d = $word.ActiveDocument

$replacements = @(
    @("291÷9=", "453÷5="),
    @("897÷6=", "297÷5="),
    @("935÷5=", "559÷6="),
    @("122÷3=", "682÷2="),
    @("643÷8=", "216÷4="),
    @("324÷8=", "359÷9="),
    @("208÷6=", "719÷5="),
    @("841÷3=", "411÷4="),
    @("646÷2=", "237÷6="),
    @("475÷6=", "925÷6="),
    @("162÷7=", "656÷3="),
    @("874÷9=", "293÷6="),
    @("472÷8=", "761÷3="),
    @("140÷2=", "817÷7="),
    @("729÷4=", "344÷2="),
    @("157÷3=", "776÷6="),
    @("938÷2=", "778÷3="),
    @("611÷3=", "163÷3="),
    @("906÷3=", "405÷3="),
    @("696÷5=", "262÷6="),
    @("130÷8=", "955÷8="),
    @("695÷4=", "803÷7="),
    @("847÷5=", "955÷9="),
    @("826÷8=", "912÷6="),
    @("358÷8=", "522÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
